$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.592.47"
$ws.Range("E2").Value = "  -0.99%  "
$ws.Range("D3").Value = "3.441.57"
$ws.Range("E3").Value = "  -1.29%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'573.39"
$ws.Range("E5").Value = "  -1.15%  "
$ws.Range("D6").Value = "'158.89"
$ws.Range("E6").Value = "  -2.46%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "3.439.82"
$ws.Range("E8").Value = "  -1.37%  "
$ws.Range("D9").Value = "'0.573"
$ws.Range("E9").Value = "  -7.25%  "
$ws.Range("E10").Value = "  -0.72%  "
$ws.Range("E11").Value = "  -3.58%  "
$ws.Range("D12").Value = "'0.437"
$ws.Range("E12").Value = "  -2.70%  "
$ws.Range("D13").Value = "4.033.36"
$ws.Range("E13").Value = "  -1.33%  "
$ws.Range("E14").Value = "  -0.49%  "
$ws.Range("D15").Value = "'27.44"
$ws.Range("E15").Value = "  -4.55%  "
$ws.Range("E16").Value = "  -10.28%  "
$ws.Range("D17").Value = "64.662.28"
$ws.Range("E17").Value = "  -0.91%  "
$ws.Range("D18").Value = "3.445.75"
$ws.Range("E18").Value = "  -0.97%  "
$ws.Range("D19").Value = "'6.15"
$ws.Range("E19").Value = "  -5.26%  "
$ws.Range("D20").Value = "'13.68"
$ws.Range("E20").Value = "  -5.17%  "
$ws.Range("D21").Value = "'377.85"
$ws.Range("E21").Value = "  -1.43%  "
$ws.Range("D22").Value = "'7.89"
$ws.Range("E22").Value = "  -4.11%  "
$ws.Range("E23").Value = "  -0.33%  "
$ws.Range("D24").Value = "'72.00"
$ws.Range("E24").Value = "  -0.95%  "
$ws.Range("D25").Value = "'0.529"
$ws.Range("E25").Value = "  -4.67%  "
$ws.Range("E26").Value = "  -0.99%  "
$ws.Range("D27").Value = "'9.90"
$ws.Range("E27").Value = "  -1.37%  "
$ws.Range("E28").Value = "  -0.54%  "
$ws.Range("E30").Value = "  -7.28%  "
$ws.Range("E31").Value = "  -2.59%  "
$ws.Range("E32").Value = "  -2.74%  "
$ws.Range("D33").Value = "'23.16"
$ws.Range("D34").Value = "'6.96"
$ws.Range("E34").Value = "  -3.39%  "
$ws.Range("E35").Value = "  -4.76%  "
$ws.Range("D36").Value = "'160.52"
$ws.Range("E36").Value = "  -1.04%  "
$ws.Range("E37").Value = "  -3.33%  "
$ws.Range("D38").Value = "2.885.84"
$ws.Range("E38").Value = "  -4.48%  "
$ws.Range("D39").Value = "'0.0744"
$ws.Range("E39").Value = "  -4.75%  "
$ws.Range("D40").Value = "'26.00"
$ws.Range("E40").Value = "  -3.48%  "
$ws.Range("D41").Value = "'42.91"
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("E42").Value = "  +0.61%  "
$ws.Range("E43").Value = "  -1.75%  "
$ws.Range("D44").Value = "'6.42"
$ws.Range("E44").Value = "  -6.34%  "
$ws.Range("D45").Value = "'25.80"
$ws.Range("E45").Value = "  -0.97%  "
$ws.Range("E46").Value = "  -4.35%  "
$ws.Range("E47").Value = "  +7.84%  "
$ws.Range("D48").Value = "'320.12"
$ws.Range("E48").Value = "  -0.09%  "
$ws.Range("E49").Value = "  -3.73%  "
$ws.Range("D50").Value = "'6.43"
$ws.Range("E50").Value = "  -4.38%  "
$ws.Range("D51").Value = "'0.841"
$ws.Range("E51").Value = "  -4.23%  "
